$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E2").Style = "Normal"
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "28.156.31"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D2:E2").Style = "Normal"

$ws.Range("D3:E3").Style = "Normal"
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "1.796.86"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("D3:E3").Style = "Normal"

$ws.Range("D4:E4").Style = "Normal"
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D4:E4").Style = "Normal"

$ws.Range("D5:E5").Style = "Normal"
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "314.45"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D5:E5").Style = "Normal"

$ws.Range("E6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7:E7").Style = "Normal"
$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5209"
$ws.Range("E7").Value = "  +2.03%  "
$ws.Range("D7:E7").Style = "Normal"

$ws.Range("D8:E8").Style = "Normal"
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3821"
$ws.Range("E8").Value = "  -3.39%  "
$ws.Range("D8:E8").Style = "Normal"

$ws.Range("D9:E9").Style = "Normal"
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07918"
$ws.Range("E9").Value = "  -3.98%  "
$ws.Range("D9:E9").Style = "Normal"

$ws.Range("D10:E10").Style = "Normal"
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "41.35"
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("D10:E10").Style = "Normal"

$ws.Range("D11:E11").Style = "Normal"
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "1.100"
$ws.Range("E11").Value = "  -1.18%  "
$ws.Range("D11:E11").Style = "Normal"

$ws.Range("D12:E12").Style = "Normal"
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "6.286"
$ws.Range("E12").Value = "  -1.23%  "
$ws.Range("D12:E12").Style = "Normal"

$ws.Range("E13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14:E14").Style = "Normal"
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "20.56"
$ws.Range("E14").Value = "  -2.90%  "
$ws.Range("D14:E14").Style = "Normal"

$ws.Range("D15:E15").Style = "Normal"
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "1.789.95"
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("D15:E15").Style = "Normal"

$ws.Range("D16:E16").Style = "Normal"
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "7.239"
$ws.Range("E16").Value = "  -4.17%  "
$ws.Range("D16:E16").Style = "Normal"

$ws.Range("D17:E17").Style = "Normal"
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "93.27"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D17:E17").Style = "Normal"

$ws.Range("D18:E18").Style = "Normal"
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001086"
$ws.Range("E18").Value = "  -3.56%  "
$ws.Range("D18:E18").Style = "Normal"

$ws.Range("E19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.41%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21:E21").Style = "Normal"
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "17.29"
$ws.Range("E21").Value = "  -3.07%  "
$ws.Range("D21:E21").Style = "Normal"

$ws.Range("D22:E22").Style = "Normal"
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "5.954"
$ws.Range("E22").Value = "  -2.46%  "
$ws.Range("D22:E22").Style = "Normal"

$ws.Range("D23:E23").Style = "Normal"
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "28.202.70"
$ws.Range("E23").Value = "  -1.33%  "
$ws.Range("D23:E23").Style = "Normal"

$ws.Range("D24:E24").Style = "Normal"
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "11.13"
$ws.Range("E24").Value = "  -2.63%  "
$ws.Range("D24:E24").Style = "Normal"

$ws.Range("D25:E25").Style = "Normal"
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "2.266"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D25:E25").Style = "Normal"

$ws.Range("D26:E26").Style = "Normal"
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "160.83"
$ws.Range("E26").Value = "  +2.57%  "
$ws.Range("D26:E26").Style = "Normal"

$ws.Range("D27:E27").Style = "Normal"
$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = "20.45"
$ws.Range("E27").Value = "  -4.14%  "
$ws.Range("D27:E27").Style = "Normal"

$ws.Range("D28:E28").Style = "Normal"
$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D28").Value = "1.999.94"
$ws.Range("E28").Value = "  -1.50%  "
$ws.Range("D28:E28").Style = "Normal"

$ws.Range("D29:E29").Style = "Normal"
$ws.Range("D29:E29").NumberFormat = "@"
$ws.Range("D29").Value = "2.337"
$ws.Range("E29").Value = "  -3.37%  "
$ws.Range("D29:E29").Style = "Normal"

$ws.Range("D30:E30").Style = "Normal"
$ws.Range("D30:E30").NumberFormat = "@"
$ws.Range("D30").Value = "123.10"
$ws.Range("E30").Value = "  -2.82%  "
$ws.Range("D30:E30").Style = "Normal"

$ws.Range("D31:E31").Style = "Normal"
$ws.Range("D31:E31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1066"
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D31:E31").Style = "Normal"

$ws.Range("D32:E32").Style = "Normal"
$ws.Range("D32:E32").NumberFormat = "@"
$ws.Range("D32").Value = "1.052"
$ws.Range("E32").Value = "  -5.67%  "
$ws.Range("D32:E32").Style = "Normal"

$ws.Range("E33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34:E34").Style = "Normal"
$ws.Range("D34:E34").NumberFormat = "@"
$ws.Range("D34").Value = "5.571"
$ws.Range("E34").Value = "  -3.69%  "
$ws.Range("D34:E34").Style = "Normal"

$ws.Range("D35:E35").Style = "Normal"
$ws.Range("D35:E35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07303"
$ws.Range("E35").Value = "  +3.47%  "
$ws.Range("D35:E35").Style = "Normal"

$ws.Range("D36:E36").Style = "Normal"
$ws.Range("D36:E36").NumberFormat = "@"
$ws.Range("D36").Value = "12.23"
$ws.Range("E36").Value = "  +8.43%  "
$ws.Range("D36:E36").Style = "Normal"

$ws.Range("D37:E37").Style = "Normal"
$ws.Range("D37:E37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02319"
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("D37:E37").Style = "Normal"

$ws.Range("D38:E38").Style = "Normal"
$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2142"
$ws.Range("E38").Value = "  -4.03%  "
$ws.Range("D38:E38").Style = "Normal"

$ws.Range("D39:E39").Style = "Normal"
$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = "8.677"
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("D39:E39").Style = "Normal"

$ws.Range("D40:E40").Style = "Normal"
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "5.068"
$ws.Range("E40").Value = "  -3.73%  "
$ws.Range("D40:E40").Style = "Normal"

$ws.Range("D41:E41").Style = "Normal"
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6156"
$ws.Range("E41").Value = "  -3.10%  "
$ws.Range("D41:E41").Style = "Normal"

$ws.Range("D42:E42").Style = "Normal"
$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "1.160"
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("D42:E42").Style = "Normal"

$ws.Range("D43:E43").Style = "Normal"
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "1.371"
$ws.Range("E43").Value = "  -2.42%  "
$ws.Range("D43:E43").Style = "Normal"

$ws.Range("D44:E44").Style = "Normal"
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "13.21"
$ws.Range("E44").Value = "  -2.76%  "
$ws.Range("D44:E44").Style = "Normal"

$ws.Range("D45:E45").Style = "Normal"
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "3.783"
$ws.Range("E45").Value = "  +1.25%  "
$ws.Range("D45:E45").Style = "Normal"

$ws.Range("D46:E46").Style = "Normal"
$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5974"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D46:E46").Style = "Normal"

$ws.Range("D47:E47").Style = "Normal"
$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = "127.74"
$ws.Range("E47").Value = "  +1.95%  "
$ws.Range("D47:E47").Style = "Normal"

$ws.Range("D48:E48").Style = "Normal"
$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "1.233"
$ws.Range("E48").Value = "  +3.14%  "
$ws.Range("D48:E48").Style = "Normal"

$ws.Range("D49:E49").Style = "Normal"
$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = "1.919"
$ws.Range("E49").Value = "  -3.82%  "
$ws.Range("D49:E49").Style = "Normal"

$ws.Range("D50:E50").Style = "Normal"
$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06768"
$ws.Range("E50").Value = "  -2.51%  "
$ws.Range("D50:E50").Style = "Normal"

$ws.Range("D51:E51").Style = "Normal"
$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("D51").Value = "73.25"
$ws.Range("E51").Value = "  -1.43%  "
$ws.Range("D51:E51").Style = "Normal"

